# Update the two-digit multiplication answer table with the newly
# generated values. Cells are addressed by (row, column) in the single
# table on the page so the update is unambiguous even though some of
# the new values coincide with old values used elsewhere in the table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "67×29=1943" },
    @{ Row = 1;  Col = 2; Text = "79×41=3239" },
    @{ Row = 1;  Col = 3; Text = "64×73=4672" },
    @{ Row = 1;  Col = 4; Text = "14×43=602" },
    @{ Row = 1;  Col = 5; Text = "15×95=1425" },

    @{ Row = 5;  Col = 1; Text = "28×34=952" },
    @{ Row = 5;  Col = 2; Text = "90×17=1530" },
    @{ Row = 5;  Col = 3; Text = "93×23=2139" },
    @{ Row = 5;  Col = 4; Text = "48×66=3168" },
    @{ Row = 5;  Col = 5; Text = "42×40=1680" },

    @{ Row = 10; Col = 1; Text = "78×23=1794" },
    @{ Row = 10; Col = 2; Text = "51×50=2550" },
    @{ Row = 10; Col = 3; Text = "46×13=598" },
    @{ Row = 10; Col = 4; Text = "85×35=2975" },
    @{ Row = 10; Col = 5; Text = "96×87=8352" },

    @{ Row = 15; Col = 1; Text = "26×15=390" },
    @{ Row = 15; Col = 2; Text = "95×17=1615" },
    @{ Row = 15; Col = 3; Text = "84×64=5376" },
    @{ Row = 15; Col = 4; Text = "31×19=589" },
    @{ Row = 15; Col = 5; Text = "94×18=1692" },

    @{ Row = 20; Col = 1; Text = "18×62=1116" },
    @{ Row = 20; Col = 2; Text = "12×34=408" },
    @{ Row = 20; Col = 3; Text = "12×30=360" },
    @{ Row = 20; Col = 4; Text = "86×14=1204" },
    @{ Row = 20; Col = 5; Text = "50×26=1300" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
